$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as TEXT (avoids Excel auto-converting
# numeric-looking strings like "0.999" or "227.16" into real numbers),
# then clear the temporary number-format style so no residual
# formatting/style is left behind on the cell.
function Set-TextValue([string]$addr, [string]$val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" '39.776.77'
Set-TextValue "E2" '  +2.65%  '
Set-TextValue "D3" '2.165.33'
Set-TextValue "E3" '  +3.11%  '
Set-TextValue "D4" '0.999'
Set-TextValue "E4" '  -0.13%  '
Set-TextValue "D5" '227.16'
Set-TextValue "E5" '  -0.07%  '
Set-TextValue "D6" '0.623'
Set-TextValue "E6" '  +1.21%  '
Set-TextValue "D7" '63.00'
Set-TextValue "E7" '  +1.82%  '
Set-TextValue "E8" '  -0.02%  '
Set-TextValue "D9" '0.392'
Set-TextValue "D10" '0.0843'
Set-TextValue "E10" '  +0.42%  '
Set-TextValue "E11" '  +0.12%  '
Set-TextValue "D12" '15.91'
Set-TextValue "E12" '  +0.96%  '
Set-TextValue "D13" '2.483.43'
Set-TextValue "E13" '  +2.98%  '
Set-TextValue "D14" '21.80'
Set-TextValue "E14" '  -0.73%  '
Set-TextValue "D15" '0.808'
Set-TextValue "E15" '  +1.09%  '
Set-TextValue "E16" '  +0.36%  '
Set-TextValue "D17" '2.164.40'
Set-TextValue "E17" '  +2.93%  '
Set-TextValue "D18" '39.749.34'
Set-TextValue "E18" '  +2.78%  '
Set-TextValue "D19" '71.82'
Set-TextValue "E19" '  +0.31%  '
Set-TextValue "D20" '6.03'
Set-TextValue "E20" '  +0.28%  '
Set-TextValue "D21" '0.0₃0849'
Set-TextValue "E21" '  +0.86%  '
Set-TextValue "D22" '228.31'
Set-TextValue "E22" '  +0.68%  '
Set-TextValue "E23" '  +0.08%  '
Set-TextValue "E24" '  +2.27%  '
Set-TextValue "D25" '2.32'
Set-TextValue "E25" '  -0.10%  '
Set-TextValue "D26" '171.60'
Set-TextValue "E26" '  +0.91%  '
Set-TextValue "D27" '9.44'
Set-TextValue "E27" '  -2.15%  '
Set-TextValue "E28" '  +2.60%  '
Set-TextValue "E29" '  +1.74%  '
Set-TextValue "D30" '19.66'
Set-TextValue "E30" '  +1.72%  '
Set-TextValue "D31" '2.70'
Set-TextValue "E31" '  +6.43%  '
Set-TextValue "D32" '0.122'
Set-TextValue "E32" '  +0.96%  '
Set-TextValue "E33" '  +0.78%  '
Set-TextValue "D34" '4.70'
Set-TextValue "E34" '  -2.18%  '
Set-TextValue "D35" '6.96'
Set-TextValue "E35" '  -2.59%  '
Set-TextValue "E36" '  +0.65%  '
Set-TextValue "D37" '3.73'
Set-TextValue "E37" '  +7.09%  '
Set-TextValue "E38" '  +2.18%  '
Set-TextValue "E39" '  -0.19%  '
Set-TextValue "D40" '4.86'
Set-TextValue "E40" '  +16.77%  '
Set-TextValue "E41" '  +1.26%  '
Set-TextValue "E42" '  -1.05%  '
Set-TextValue "D43" '17.57'
Set-TextValue "E43" '  -2.19%  '
Set-TextValue "D44" '1.514.10'
Set-TextValue "E44" '  -0.67%  '
Set-TextValue "E45" '  -0.18%  '
Set-TextValue "E46" '  +2.58%  '
Set-TextValue "D47" '0.0929'
Set-TextValue "E47" '  +2.05%  '
Set-TextValue "E48" '  +0.00%  '
Set-TextValue "E49" '  +1.29%  '
Set-TextValue "B50" 'MultiversX'
Set-TextValue "C50" 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue "D50" '49.65'
Set-TextValue "E50" '  +8.24%  '
Set-TextValue "B51" 'TerraClassic'
Set-TextValue "C51" 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
Set-TextValue "D51" '0.000191'
Set-TextValue "E51" '  +28.76%  '
